$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": refresh the handoff timestamps recorded for the
# most recently handed-off file (62938945-5bda-4e7c-93f9-3152292fd58f).

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for that file's row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-09-03 06:44:36"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for that file's row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-09-03 06:44:32"

# de-de sheet: "Latest Handoff Datetime" column (H) for that file's row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-09-03 06:44:36"
